# Alarm Normal load method changes
# Adds a new "MinBatteryLoadingDetail" column (K) to the "Add Panels" sheet:
#   K7 -> header "MinBatteryLoadingDetail" (same style as the other header cells A7:D7)
#   K8 -> value  "Minimum battery (Ah)"   (same style as the other data cells in row 8)
# and leaves the active selection on K12, matching the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new cells. K8 is written first so that the shared-string table
# ends up with "Minimum battery (Ah)" before "MinBatteryLoadingDetail", matching
# the order they were authored in the workbook.
$ws.Range("K8").Value = "Minimum battery (Ah)"
$ws.Range("K7").Value = "MinBatteryLoadingDetail"

# Match formatting of the neighboring header/data cells.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("K7").PasteSpecial(-4122) | Out-Null

$ws.Range("B8").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null

# Restore the active cell selection as recorded in the workbook.
$ws.Range("K12").Select() | Out-Null
